# Update the localization status for 228b4934-faed-44a0-b362-1b99a5cea0b9.md
# from "Ready for handoff" to "In Translation" across the Overview,
# zh-cn and de-de report sheets (report regenerated for archive).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2:C2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "In Translation"
